$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Issue number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Type-change cells: convert between numeric and text ("N/A") representations ---
# Donor cells (unaffected by this edit) used to carry over the correct style/shared-string type:
#   C23  -> style s13, text "0"      (numeric placeholder for N/A counts)
#   E23  -> style s13, text "***.*"  (placeholder for N/A percentages)
#   G19  -> style s14 (plain integer format)
#   M20  -> style s15 (percent-change number format)

$ws.Range("C23").Copy($ws.Range("C16"))
$ws.Range("C23").Copy($ws.Range("C18"))
$ws.Range("C23").Copy($ws.Range("D20"))

$ws.Range("E23").Copy($ws.Range("E20"))

# text -> number (plain integer, style s14)
$ws.Range("G19").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("G19").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("G19").Copy($ws.Range("F27"))
$ws.Range("F27").Value = 2
$ws.Range("G19").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 3
$ws.Range("G19").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 2
$ws.Range("G19").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 2
$ws.Range("G19").Copy($ws.Range("I29"))
$ws.Range("I29").Value = 2
$ws.Range("G19").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("G19").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 1
$ws.Range("G19").Copy($ws.Range("I30"))
$ws.Range("I30").Value = 1

# text -> number (percent-change format, style s15)
$ws.Range("M20").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 100
$ws.Range("M20").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -66.666666666666

# --- Simple same-type numeric value updates ---
$ws.Range("M15").Value = -16.666666666666
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -41.666666666666
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 81
$ws.Range("K16").Value = -34.567901234567
$ws.Range("L16").Value = -41.758241758241
$ws.Range("M16").Value = -18.461538461538
$ws.Range("N16").Value = -87.410926365795
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 52
$ws.Range("J17").Value = 68
$ws.Range("K17").Value = -23.529411764705
$ws.Range("L17").Value = -51.401869158878
$ws.Range("M17").Value = 10.63829787234
$ws.Range("N17").Value = -59.375
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -65
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 135
$ws.Range("K18").Value = -44.444444444444
$ws.Range("L18").Value = -50.980392156862
$ws.Range("M18").Value = -20.212765957446
$ws.Range("N18").Value = -80.719794344473
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = 21.739130434782
$ws.Range("F19").Value = 81
$ws.Range("H19").Value = -3.571428571428
$ws.Range("I19").Value = 453
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = -9.4
$ws.Range("L19").Value = -25.859247135842
$ws.Range("M19").Value = -10.119047619047
$ws.Range("N19").Value = -61.015490533562
$ws.Range("L20").Value = -64.705882352941
$ws.Range("N20").Value = -98.382749326145
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 6.451612903225
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -18.461538461538
$ws.Range("I21").Value = 645
$ws.Range("J21").Value = 806
$ws.Range("K21").Value = -19.975186104218
$ws.Range("L21").Value = -34.451219512195
$ws.Range("M21").Value = -12.364130434782
$ws.Range("N21").Value = -73.960436011304
$ws.Range("C22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = 36.363636363636
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = -13.385826771653
$ws.Range("I24").Value = 708
$ws.Range("J24").Value = 828
$ws.Range("K24").Value = -14.492753623188
$ws.Range("L24").Value = -27.16049382716
$ws.Range("M24").Value = -8.880308880308
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = -24.137931034482
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 108
$ws.Range("H25").Value = -20.37037037037
$ws.Range("I25").Value = 516
$ws.Range("J25").Value = 683
$ws.Range("K25").Value = -24.450951683748
$ws.Range("L25").Value = -28.531855955678
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 62.5
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 6.25
$ws.Range("I26").Value = 181
$ws.Range("J26").Value = 178
$ws.Range("K26").Value = 1.685393258426
$ws.Range("L26").Value = -19.911504424778
$ws.Range("M26").Value = 50.833333333333
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = 33.333333333333
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = -14.285714285714
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 5.714285714285
$ws.Range("L28").Value = 5.714285714285
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -71.428571428571
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -85.714285714285
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -20
$ws.Range("I31").Value = 10
$ws.Range("J31").Value = 17
$ws.Range("K31").Value = -41.176470588235
$ws.Range("L31").Value = 66.666666666666
